$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Fecha (D) / Volumen (M) / Precio minimo (N) / Precio maximo (O) /
# Precio promedio ponderado (P) / Precio $/Kg (S) for rows 2-9 - the sheet
# is being re-sorted by ascending Fecha after a new week's row is added.
$rowData = @{
    2  = @{ D = 44490; M = 400; N = 9500;  O = 10000; P = 9750;  S = 4875 }
    3  = @{ D = 44461; M = 200; N = 11000; O = 12000; P = 11500; S = 5750 }
    4  = @{ D = 44455; M = 200; N = 12000; O = 13000; P = 12500; S = 6250 }
    5  = @{ D = 44489; M = 160; N = 9500;  O = 10000; P = 9750;  S = 4875 }
    6  = @{ D = 44497; M = 500; N = 9000;  O = 10000; P = 9500;  S = 4750 }
    7  = @{ D = 44482; M = 240; N = 10000; O = 11000; P = 10500; S = 5250 }
    8  = @{ D = 44475; M = 240; N = 11000; O = 12000; P = 11500; S = 5750 }
    9  = @{ D = 44454; M = 160; N = 12000; O = 13000; P = 12500; S = 6250 }
}

foreach ($r in $rowData.Keys) {
    $vals = $rowData[$r]
    $ws.Cells.Item($r, 4).Value  = $vals.D   # D - Fecha
    $ws.Cells.Item($r, 13).Value = $vals.M   # M - Volumen
    $ws.Cells.Item($r, 14).Value = $vals.N   # N - Precio minimo
    $ws.Cells.Item($r, 15).Value = $vals.O   # O - Precio maximo
    $ws.Cells.Item($r, 16).Value = $vals.P   # P - Precio promedio ponderado
    $ws.Cells.Item($r, 19).Value = $vals.S   # S - Precio $/Kg
}

# Append the new week's row (row 10). Match D's date number format so it
# renders the same as the rest of the date column.
$ws.Cells.Item(10, 4).NumberFormat = $ws.Cells.Item(9, 4).NumberFormat

$ws.Cells.Item(10, 1).Value  = 2
$ws.Cells.Item(10, 2).Value  = "Comercializadora del Agro de Limarí"
$ws.Cells.Item(10, 3).Value  = "Coquimbo"
$ws.Cells.Item(10, 4).Value  = 44517
$ws.Cells.Item(10, 5).Value  = 4
$ws.Cells.Item(10, 6).Value  = "Fruta"
$ws.Cells.Item(10, 7).Value  = 100101
$ws.Cells.Item(10, 8).Value  = "Berries"
$ws.Cells.Item(10, 9).Value  = 100101001
$ws.Cells.Item(10, 10).Value = "Arándano (blue)"
$ws.Cells.Item(10, 11).Value = "Sin especificar"
$ws.Cells.Item(10, 12).Value = "Primera"
$ws.Cells.Item(10, 13).Value = 400
$ws.Cells.Item(10, 14).Value = 5500
$ws.Cells.Item(10, 15).Value = 6000
$ws.Cells.Item(10, 16).Value = 5750
$ws.Cells.Item(10, 17).Value = "$/bandeja 2 kilos"
$ws.Cells.Item(10, 18).Value = "Provincia de Limarí"
$ws.Cells.Item(10, 19).Value = 2875
$ws.Cells.Item(10, 20).Value = 2
